$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price, 1h volume %, hour) scraped on 2023-01-21.
# Text-like columns (B/C) are plain strings; numeric-looking columns (D/E/G) are
# forced to text with a leading apostrophe (matching the original inlineStr cells),
# then the resulting quote-prefix style is reset back to Normal to avoid changing
# cell formatting/style ids.

# Row 2
$ws.Range("D2").Value = "'306.24"
$ws.Range("E2").Value = "'6.33%"
$ws.Range("G2").Value = "'10"
$ws.Range("D2:G2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'35.05"
$ws.Range("E3").Value = "'12.79%"
$ws.Range("G3").Value = "'10"
$ws.Range("D3:G3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.185"
$ws.Range("E4").Value = "'5.28%"
$ws.Range("G4").Value = "'10"
$ws.Range("D4:G4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.07879"
$ws.Range("E5").Value = "'7.49%"
$ws.Range("G5").Value = "'10"
$ws.Range("D5:G5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'2.364"
$ws.Range("E6").Value = "'6.43%"
$ws.Range("G6").Value = "'10"
$ws.Range("D6:G6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'7.991"
$ws.Range("E7").Value = "'3.23%"
$ws.Range("G7").Value = "'10"
$ws.Range("D7:G7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = "'3.953"
$ws.Range("E8").Value = "'6.12%"
$ws.Range("G8").Value = "'10"
$ws.Range("D8:G8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.9312"
$ws.Range("E9").Value = "'3.23%"
$ws.Range("G9").Value = "'10"
$ws.Range("D9:G9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.1019"
$ws.Range("E10").Value = "'11.28%"
$ws.Range("G10").Value = "'10"
$ws.Range("D10:G10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1835"
$ws.Range("E11").Value = "'8.43%"
$ws.Range("G11").Value = "'10"
$ws.Range("D11:G11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.08683"
$ws.Range("E12").Value = "'6.47%"
$ws.Range("G12").Value = "'10"
$ws.Range("D12:G12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03356"
$ws.Range("E13").Value = "'7.43%"
$ws.Range("G13").Value = "'10"
$ws.Range("D13:G13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09894"
$ws.Range("E14").Value = "'-0.46%"
$ws.Range("G14").Value = "'10"
$ws.Range("D14:G14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001478"
$ws.Range("E15").Value = "'-1.25%"
$ws.Range("G15").Value = "'10"
$ws.Range("D15:G15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.005606"
$ws.Range("E16").Value = "'-1.89%"
$ws.Range("G16").Value = "'10"
$ws.Range("D16:G16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'-1.38%"
$ws.Range("G17").Value = "'10"
$ws.Range("D17:G17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.142"
$ws.Range("E18").Value = "'2.91%"
$ws.Range("G18").Value = "'10"
$ws.Range("D18:G18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.3395"
$ws.Range("E19").Value = "'2.03%"
$ws.Range("G19").Value = "'10"
$ws.Range("D19:G19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.1302"
$ws.Range("G20").Value = "'10"
$ws.Range("D20:G20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'4.578"
$ws.Range("E21").Value = "'9.37%"
$ws.Range("G21").Value = "'10"
$ws.Range("D21:G21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.2290"
$ws.Range("E22").Value = "'9.01%"
$ws.Range("G22").Value = "'10"
$ws.Range("D22:G22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04566"
$ws.Range("E23").Value = "'0.75%"
$ws.Range("G23").Value = "'10"
$ws.Range("D23:G23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.001210"
$ws.Range("E24").Value = "'-0.02%"
$ws.Range("G24").Value = "'10"
$ws.Range("D24:G24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.004464"
$ws.Range("E25").Value = "'7.31%"
$ws.Range("G25").Value = "'10"
$ws.Range("D25:G25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.0001293"
$ws.Range("E26").Value = "'-0.63%"
$ws.Range("G26").Value = "'10"
$ws.Range("D26:G26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.0003383"
$ws.Range("E27").Value = "'-0.38%"
$ws.Range("G27").Value = "'10"
$ws.Range("D27:G27").Style = "Normal"

# Row 28
$ws.Range("G28").Value = "'10"
$ws.Range("D28:G28").Style = "Normal"

# Row 29
$ws.Range("G29").Value = "'10"
$ws.Range("D29:G29").Style = "Normal"

# Row 30
$ws.Range("G30").Value = "'10"
$ws.Range("D30:G30").Style = "Normal"

# Row 31
$ws.Range("G31").Value = "'10"
$ws.Range("D31:G31").Style = "Normal"

# Row 32
$ws.Range("G32").Value = "'10"
$ws.Range("D32:G32").Style = "Normal"

# Row 33
$ws.Range("G33").Value = "'10"
$ws.Range("D33:G33").Style = "Normal"

# Row 34
$ws.Range("G34").Value = "'10"
$ws.Range("D34:G34").Style = "Normal"

# Row 35
$ws.Range("G35").Value = "'10"
$ws.Range("D35:G35").Style = "Normal"

# Row 36
$ws.Range("G36").Value = "'10"
$ws.Range("D36:G36").Style = "Normal"

# Row 37
$ws.Range("G37").Value = "'10"
$ws.Range("D37:G37").Style = "Normal"

# Row 38
$ws.Range("G38").Value = "'10"
$ws.Range("D38:G38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.01789"
$ws.Range("E39").Value = "'13.53%"
$ws.Range("G39").Value = "'10"
$ws.Range("D39:G39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.04788"
$ws.Range("E40").Value = "'7.70%"
$ws.Range("G40").Value = "'10"
$ws.Range("D40:G40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007755"
$ws.Range("E41").Value = "'5.75%"
$ws.Range("G41").Value = "'10"
$ws.Range("D41:G41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1415"
$ws.Range("E42").Value = "'6.46%"
$ws.Range("G42").Value = "'10"
$ws.Range("D42:G42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.007050"
$ws.Range("G43").Value = "'10"
$ws.Range("D43:G43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.002205"
$ws.Range("E44").Value = "'-3.79%"
$ws.Range("G44").Value = "'10"
$ws.Range("D44:G44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.009500"
$ws.Range("E45").Value = "'14.19%"
$ws.Range("G45").Value = "'10"
$ws.Range("D45:G45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00005967"
$ws.Range("E46").Value = "'-2.20%"
$ws.Range("G46").Value = "'10"
$ws.Range("D46:G46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.00000000744"
$ws.Range("E47").Value = "'-0.87%"
$ws.Range("G47").Value = "'10"
$ws.Range("D47:G47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'10.71%"
$ws.Range("G48").Value = "'10"
$ws.Range("D48:G48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.002676"
$ws.Range("E49").Value = "'33.69%"
$ws.Range("G49").Value = "'10"
$ws.Range("D49:G49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.00002083"
$ws.Range("E50").Value = "'-0.87%"
$ws.Range("G50").Value = "'10"
$ws.Range("D50:G50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0001984"
$ws.Range("E51").Value = "'-0.87%"
$ws.Range("G51").Value = "'10"
$ws.Range("D51:G51").Style = "Normal"
